$d = $word.ActiveDocument

$replacements = @(
    @("2025-03-25 Tuesday", "2025-03-26 Wednesday"),
    @("229×3=687", "786×3=2358"),
    @("401×4=1604", "610×3=1830"),
    @("490×8=3920", "786×4=3144"),
    @("419×6=2514", "830×9=7470"),
    @("716×4=2864", "224×5=1120"),
    @("178×9=1602", "465×7=3255"),
    @("345×7=2415", "856×4=3424"),
    @("486×8=3888", "791×2=1582"),
    @("226×3=678", "256×5=1280"),
    @("843×3=2529", "300×4=1200"),
    @("971×8=7768", "722×4=2888"),
    @("931×9=8379", "607×8=4856"),
    @("913×3=2739", "910×7=6370"),
    @("925×2=1850", "348×5=1740"),
    @("284×3=852", "729×5=3645"),
    @("623×8=4984", "590×8=4720"),
    @("803×9=7227", "253×5=1265"),
    @("920×7=6440", "409×8=3272"),
    @("540×3=1620", "227×2=454"),
    @("820×2=1640", "735×9=6615"),
    @("386×4=1544", "945×2=1890"),
    @("878×5=4390", "632×8=5056"),
    @("514×8=4112", "603×3=1809"),
    @("424×7=2968", "552×2=1104"),
    @("216×9=1944", "737×9=6633")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
